$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Simulations S2.1 (own RES inve)")
$ws2 = $wb.Worksheets.Item("Simulations S2.2 (PPAs)")

# ---------------------------------------------------------------------------
# Sheet "Simulations S2.1 (own RES inve)" - BoP capex/opex formula updates
# ---------------------------------------------------------------------------

# Row 5 - BoP capex: switch E (80%) column formula to a direct 75%-of-B formula
$ws1.Range("E5").Formula = "=B5*0.75"

# Row 7 - de-share the existing 80%-of-B formula (same result, own formula)
$ws1.Range("E7").Formula = "=B7*0.8"

# Row 9 - de-share D (90%-of-B, same result) and switch E to 75%-of-B
$ws1.Range("D9").Formula = "=B9*0.9"
$ws1.Range("E9").Formula = "=B9*0.75"

# Row 21 - E now 65% of B instead of 80% of D
$ws1.Range("E21").Formula = "=B21*0.65"

# Row 22 - de-share D (90%-of-B, same result) and switch E to 65% of B
$ws1.Range("D22").Formula = "=B22*0.9"
$ws1.Range("E22").Formula = "=B22*0.65"

# ---------------------------------------------------------------------------
# Sheet "Simulations S2.2 (PPAs)" - same BoP capex/opex formula updates
# ---------------------------------------------------------------------------

$ws2.Range("E4").Formula = "=B4*0.75"
$ws2.Range("E8").Formula = "=B8*0.75"
$ws2.Range("E21").Formula = "=B21*0.65"
$ws2.Range("E22").Formula = "=B22*0.65"

# Row 23 - E23 goes back to the plain "=$B$23" style value/format used by C23/D23
$ws2.Range("C23").Copy()
$ws2.Range("E23").PasteSpecial(-4122)
$ws2.Range("E23").Value = 1.33

# ---------------------------------------------------------------------------
# View state - active sheet / selection moves from S2.2 to S2.1
# ---------------------------------------------------------------------------

[void]$ws2.Range("E9").Select()
[void]$ws1.Activate()
[void]$ws1.Range("E8").Select()
